# The "Free Food" question (row 74 - "Do you currently have access to free
# food for your household?...") was removed from the question bank sheet.
# Deleting the whole row shifts every subsequent row up by one and Excel
# automatically cleans up the now-unused shared strings for that row's three
# cells, renumbering the sharedStrings table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("74:74").Delete()

# Reflect the new selection/view state left behind after the edit.
$ws.Range("A74:XFD80").Select()
